$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.067.14"
$ws.Range("E2").Value = "  +9.44%  "

$ws.Range("D3").Value = "'2.697.64"
$ws.Range("E3").Value = "  +11.14%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'187.78"
$ws.Range("E5").Value = "  +12.75%  "

$ws.Range("D6").Value = "'588.75"
$ws.Range("E6").Value = "  +4.28%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'0.541"
$ws.Range("E8").Value = "  +5.11%  "

$ws.Range("D9").Value = "'0.196"
$ws.Range("E9").Value = "  +15.74%  "

$ws.Range("D10").Value = "'2.695.99"
$ws.Range("E10").Value = "  +11.10%  "

$ws.Range("E11").Value = "  +1.47%  "

$ws.Range("D12").Value = "'0.358"
$ws.Range("E12").Value = "  +6.94%  "

$ws.Range("D13").Value = "'4.74"
$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("D14").Value = "'75.774.04"
$ws.Range("E14").Value = "  +9.48%  "

$ws.Range("D15").Value = "'3.190.28"
$ws.Range("E15").Value = "  +11.13%  "

$ws.Range("D16").Value = "'0.0000189"
$ws.Range("E16").Value = "  +6.52%  "

$ws.Range("D17").Value = "'26.57"
$ws.Range("E17").Value = "  +10.85%  "

$ws.Range("D18").Value = "'2.702.31"
$ws.Range("E18").Value = "  +11.81%  "

$ws.Range("D19").Value = "'9.36"
$ws.Range("E19").Value = "  +31.09%  "

$ws.Range("D20").Value = "'12.00"
$ws.Range("E20").Value = "  +11.06%  "

$ws.Range("D21").Value = "'374.92"
$ws.Range("E21").Value = "  +9.50%  "

$ws.Range("D22").Value = "'2.28"
$ws.Range("E22").Value = "  +16.04%  "

$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = "  +4.82%  "

$ws.Range("D24").Value = "'6.28"
$ws.Range("E24").Value = "  +4.65%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "'70.06"
$ws.Range("E26").Value = "  +6.01%  "

$ws.Range("D27").Value = "'4.19"
$ws.Range("E27").Value = "  +9.89%  "

$ws.Range("E28").Value = "  +10.23%  "

$ws.Range("D29").Value = "'2.825.59"
$ws.Range("E29").Value = "  +10.60%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("D31").Value = "'0.0₃0952"
$ws.Range("E31").Value = "  +12.01%  "

$ws.Range("D32").Value = "'522.48"
$ws.Range("E32").Value = "  +15.17%  "

$ws.Range("D33").Value = "'1.42"
$ws.Range("E33").Value = "  +14.29%  "

$ws.Range("D34").Value = "'7.79"
$ws.Range("E34").Value = "  +5.18%  "

$ws.Range("D35").Value = "'1.75"
$ws.Range("E35").Value = "  +8.80%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.120"
$ws.Range("E37").Value = "  +8.71%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'162.74"
$ws.Range("E38").Value = "  +2.47%  "

$ws.Range("D39").Value = "'19.29"
$ws.Range("E39").Value = "  +5.89%  "

$ws.Range("D40").Value = "'19.41"
$ws.Range("E40").Value = "  +1.56%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "'5.04"
$ws.Range("E42").Value = "  +14.47%  "

$ws.Range("D43").Value = "'170.25"
$ws.Range("E43").Value = "  +26.10%  "

$ws.Range("E44").Value = "  +11.98%  "

$ws.Range("E45").Value = "  +9.54%  "

$ws.Range("D46").Value = "'1.19"
$ws.Range("E46").Value = "  +9.33%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.38"
$ws.Range("E47").Value = "  +14.01%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'39.33"
$ws.Range("E48").Value = "  +3.85%  "

$ws.Range("D49").Value = "'0.0845"
$ws.Range("E49").Value = "  +16.66%  "

$ws.Range("D50").Value = "'3.66"
$ws.Range("E50").Value = "  +7.66%  "

$ws.Range("D51").Value = "'0.539"
$ws.Range("E51").Value = "  +10.51%  "
